$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift existing week columns (B:D) three columns to the right ---
# Before: B=Jun_17, C=Jun_15, D=Jun_13, E=UN/Jun_10  (data rows hold "UN" or rating-change text)
# After we want the old B:E block living at E:H, freeing up B:D for two new report weeks
# (Jun_27 / Jun_26 / Jun_26).
$ws.Range("B:D").Insert()

# --- 2. New header row values for the freshly inserted week columns ---
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- 3. Fill the new B:D columns with "UN" (unchanged) for every analyst row ---
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# --- 4. Row 22 (BidaskClub) records a real downgrade this period ---
$downgradeText = "6/23/2018,Downgrades,Strong-Buy -> Buy,"
$ws.Range("C22").Value = $downgradeText
$ws.Range("D22").Value = $downgradeText
$ws.Range("D22").Interior.ColorIndex = 45

# --- 5. New analyst/firm rows appended at the bottom ---
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# --- 6. Column widths: C:H all get the same fixed "8" width used by the report ---
$ws.Columns.Item(3).ColumnWidth = 7.14
$ws.Columns.Item(4).ColumnWidth = 7.14
$ws.Columns.Item(5).ColumnWidth = 7.14
$ws.Columns.Item(6).ColumnWidth = 7.14
$ws.Columns.Item(7).ColumnWidth = 7.14
$ws.Columns.Item(8).ColumnWidth = 7.14
